$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 167624
$ws.Cells.Item(2, 5).Value = 14911
$ws.Cells.Item(2, 6).Value = 14911
$ws.Cells.Item(2, 7).Value = 10988
$ws.Cells.Item(2, 8).Value = 7823
$ws.Cells.Item(2, 9).Value = 7648
$ws.Cells.Item(2, 10).Value = 175
$ws.Cells.Item(2, 11).Value = 289338
$ws.Cells.Item(2, 12).Value = 151190
$ws.Cells.Item(2, 13).Value = 138148
$ws.Cells.Item(2, 14).Value = 136239
$ws.Cells.Item(2, 15).Value = 1909
$ws.Cells.Item(2, 16).Value = 5827
$ws.Cells.Item(2, 17).Value = 19366
$ws.Cells.Item(2, 18).Value = -11314
$ws.Cells.Item(2, 19).Value = -8754
$ws.Cells.Item(2, 20).Value = 10538
$ws.Cells.Item(2, 21).Value = 8828
$ws.Cells.Item(2, 22).Value = 122496
$ws.Cells.Item(2, 23).Value = 8.9
$ws.Cells.Item(2, 24).Value = 4.67
$ws.Cells.Item(2, 25).Value = 5.7
$ws.Cells.Item(2, 26).Value = 2.68
$ws.Cells.Item(2, 27).Value = 109.44
$ws.Cells.Item(2, 28).Value = 2069.77
$ws.Cells.Item(2, 29).Value = 6562
$ws.Cells.Item(2, 30).Value = 9.68
$ws.Cells.Item(2, 31).Value = 117571
$ws.Cells.Item(2, 32).Value = 0.54
$ws.Cells.Item(2, 33).Value = 750
$ws.Cells.Item(2, 34).Value = 1.18
$ws.Cells.Item(2, 35).Value = 11.36
$ws.Cells.Item(2, 36).Value = 116549784

# Row 3
$ws.Cells.Item(3, 4).Value = 161325
$ws.Cells.Item(3, 5).Value = 14641
$ws.Cells.Item(3, 6).Value = 14641
$ws.Cells.Item(3, 7).Value = 9267
$ws.Cells.Item(3, 8).Value = 7392
$ws.Cells.Item(3, 9).Value = 7336
$ws.Cells.Item(3, 10).Value = 55
$ws.Cells.Item(3, 11).Value = 319364
$ws.Cells.Item(3, 12).Value = 164164
$ws.Cells.Item(3, 13).Value = 155200
$ws.Cells.Item(3, 14).Value = 152420
$ws.Cells.Item(3, 15).Value = 2780
$ws.Cells.Item(3, 16).Value = 6672
$ws.Cells.Item(3, 17).Value = 30759
$ws.Cells.Item(3, 18).Value = -22050
$ws.Cells.Item(3, 19).Value = -7462
$ws.Cells.Item(3, 20).Value = 22167
$ws.Cells.Item(3, 21).Value = 8592
$ws.Cells.Item(3, 22).Value = 129363
$ws.Cells.Item(3, 23).Value = 9.08
$ws.Cells.Item(3, 24).Value = 4.58
$ws.Cells.Item(3, 25).Value = 5.08
$ws.Cells.Item(3, 26).Value = 2.43
$ws.Cells.Item(3, 27).Value = 105.78
$ws.Cells.Item(3, 28).Value = 2055.8
$ws.Cells.Item(3, 29).Value = 5866
$ws.Cells.Item(3, 30).Value = 8.52
$ws.Cells.Item(3, 31).Value = 116920
$ws.Cells.Item(3, 32).Value = 0.43
$ws.Cells.Item(3, 33).Value = 750
$ws.Cells.Item(3, 34).Value = 1.5
$ws.Cells.Item(3, 35).Value = 13.33
$ws.Cells.Item(3, 36).Value = 133445785

# Row 4
$ws.Cells.Item(4, 4).Value = 166915
$ws.Cells.Item(4, 5).Value = 14454
$ws.Cells.Item(4, 6).Value = 14450
$ws.Cells.Item(4, 7).Value = 11740
$ws.Cells.Item(4, 8).Value = 8671
$ws.Cells.Item(4, 9).Value = 8502
$ws.Cells.Item(4, 10).Value = 170
$ws.Cells.Item(4, 11).Value = 325408
$ws.Cells.Item(4, 12).Value = 160560
$ws.Cells.Item(4, 13).Value = 164848
$ws.Cells.Item(4, 14).Value = 161762
$ws.Cells.Item(4, 15).Value = 3085
$ws.Cells.Item(4, 16).Value = 6672
$ws.Cells.Item(4, 17).Value = 29166
$ws.Cells.Item(4, 18).Value = -19758
$ws.Cells.Item(4, 19).Value = -10169
$ws.Cells.Item(4, 20).Value = 19757
$ws.Cells.Item(4, 21).Value = 9409
$ws.Cells.Item(4, 22).Value = 120718
$ws.Cells.Item(4, 23).Value = 8.66
$ws.Cells.Item(4, 24).Value = 5.2
$ws.Cells.Item(4, 25).Value = 5.41
$ws.Cells.Item(4, 26).Value = 2.69
$ws.Cells.Item(4, 27).Value = 97.40000000000001
$ws.Cells.Item(4, 28).Value = 2183.46
$ws.Cells.Item(4, 29).Value = 6371
$ws.Cells.Item(4, 30).Value = 8.949999999999999
$ws.Cells.Item(4, 31).Value = 123228
$ws.Cells.Item(4, 32).Value = 0.46
$ws.Cells.Item(4, 33).Value = 750
$ws.Cells.Item(4, 34).Value = 1.32
$ws.Cells.Item(4, 35).Value = 11.58
$ws.Cells.Item(4, 36).Value = 133445785

# Row 5
$ws.Cells.Item(5, 4).Value = 191660
$ws.Cells.Item(5, 5).Value = 13676
$ws.Cells.Item(5, 6).Value = 13676
$ws.Cells.Item(5, 7).Value = 10813
$ws.Cells.Item(5, 8).Value = 7275
$ws.Cells.Item(5, 9).Value = 7161
$ws.Cells.Item(5, 10).Value = 114
$ws.Cells.Item(5, 11).Value = 333738
$ws.Cells.Item(5, 12).Value = 163383
$ws.Cells.Item(5, 13).Value = 170355
$ws.Cells.Item(5, 14).Value = 167155
$ws.Cells.Item(5, 15).Value = 3200
$ws.Cells.Item(5, 16).Value = 6672
$ws.Cells.Item(5, 17).Value = 17198
$ws.Cells.Item(5, 18).Value = -12342
$ws.Cells.Item(5, 19).Value = -4503
$ws.Cells.Item(5, 20).Value = 12011
$ws.Cells.Item(5, 21).Value = 5187
$ws.Cells.Item(5, 22).Value = 115492
$ws.Cells.Item(5, 23).Value = 7.13
$ws.Cells.Item(5, 24).Value = 3.8
$ws.Cells.Item(5, 25).Value = 4.35
$ws.Cells.Item(5, 26).Value = 2.21
$ws.Cells.Item(5, 27).Value = 95.91
$ws.Cells.Item(5, 28).Value = 2276.91
$ws.Cells.Item(5, 29).Value = 5366
$ws.Cells.Item(5, 30).Value = 10.92
$ws.Cells.Item(5, 31).Value = 127336
$ws.Cells.Item(5, 32).Value = 0.46
$ws.Cells.Item(5, 33).Value = 750
$ws.Cells.Item(5, 34).Value = 1.28
$ws.Cells.Item(5, 35).Value = 13.75
$ws.Cells.Item(5, 36).Value = 133445785

# Row 6
$ws.Cells.Item(6, 4).Value = 207804
$ws.Cells.Item(6, 5).Value = 10261
$ws.Cells.Item(6, 6).Value = 10261
$ws.Cells.Item(6, 7).Value = 5698
$ws.Cells.Item(6, 8).Value = 4080
$ws.Cells.Item(6, 9).Value = 3987
$ws.Cells.Item(6, 11).Value = 333049
$ws.Cells.Item(6, 12).Value = 162818
$ws.Cells.Item(6, 13).Value = 170231
$ws.Cells.Item(6, 14).Value = 166960
$ws.Cells.Item(6, 16).Value = 6672
$ws.Cells.Item(6, 17).Value = 15753
$ws.Cells.Item(6, 18).Value = -11648
$ws.Cells.Item(6, 19).Value = -4185
$ws.Cells.Item(6, 20).Value = 11980
$ws.Cells.Item(6, 21).Value = 3773
$ws.Cells.Item(6, 22).Value = 112858
$ws.Cells.Item(6, 23).Value = 4.94
$ws.Cells.Item(6, 24).Value = 1.96
$ws.Cells.Item(6, 25).Value = 2.39
$ws.Cells.Item(6, 26).Value = 1.22
$ws.Cells.Item(6, 27).Value = 95.64
$ws.Cells.Item(6, 28).Value = 2319.31
$ws.Cells.Item(6, 29).Value = 2988
$ws.Cells.Item(6, 30).Value = 15.15
$ws.Cells.Item(6, 31).Value = 127188
$ws.Cells.Item(6, 32).Value = 0.36
$ws.Cells.Item(6, 33).Value = 750
$ws.Cells.Item(6, 34).Value = 1.66
$ws.Cells.Item(6, 35).Value = 24.7
$ws.Cells.Item(6, 36).Value = 133445785

# Row 7
$ws.Cells.Item(7, 4).Value = 206333
$ws.Cells.Item(7, 5).Value = 4987
$ws.Cells.Item(7, 7).Value = 1704
$ws.Cells.Item(7, 8).Value = 695
$ws.Cells.Item(7, 9).Value = 636
$ws.Cells.Item(7, 11).Value = 337241
$ws.Cells.Item(7, 12).Value = 165458
$ws.Cells.Item(7, 13).Value = 171783
$ws.Cells.Item(7, 14).Value = 168436
$ws.Cells.Item(7, 16).Value = 6670
$ws.Cells.Item(7, 17).Value = 16268
$ws.Cells.Item(7, 18).Value = -11113
$ws.Cells.Item(7, 19).Value = -539
$ws.Cells.Item(7, 20).Value = 11228
$ws.Cells.Item(7, 21).Value = 6267
$ws.Cells.Item(7, 23).Value = 2.42
$ws.Cells.Item(7, 24).Value = 0.34
$ws.Cells.Item(7, 25).Value = 0.38
$ws.Cells.Item(7, 26).Value = 0.21
$ws.Cells.Item(7, 27).Value = 96.31999999999999
$ws.Cells.Item(7, 29).Value = 477
$ws.Cells.Item(7, 30).Value = 61.97
$ws.Cells.Item(7, 31).Value = 128312
$ws.Cells.Item(7, 32).Value = 0.23
$ws.Cells.Item(7, 33).Value = 734
$ws.Cells.Item(7, 34).Value = 2.49
$ws.Cells.Item(7, 35).Value = 154.01

# Row 8
$ws.Cells.Item(8, 4).Value = 201890
$ws.Cells.Item(8, 5).Value = 5816
$ws.Cells.Item(8, 7).Value = 3126
$ws.Cells.Item(8, 8).Value = 2295
$ws.Cells.Item(8, 9).Value = 2258
$ws.Cells.Item(8, 11).Value = 334992
$ws.Cells.Item(8, 12).Value = 161817
$ws.Cells.Item(8, 13).Value = 173175
$ws.Cells.Item(8, 14).Value = 169824
$ws.Cells.Item(8, 16).Value = 6670
$ws.Cells.Item(8, 17).Value = 18660
$ws.Cells.Item(8, 18).Value = -11993
$ws.Cells.Item(8, 19).Value = -4501
$ws.Cells.Item(8, 20).Value = 11400
$ws.Cells.Item(8, 21).Value = 6752
$ws.Cells.Item(8, 23).Value = 2.88
$ws.Cells.Item(8, 24).Value = 1.14
$ws.Cells.Item(8, 25).Value = 1.34
$ws.Cells.Item(8, 26).Value = 0.68
$ws.Cells.Item(8, 27).Value = 93.44
$ws.Cells.Item(8, 29).Value = 1692
$ws.Cells.Item(8, 30).Value = 16.73
$ws.Cells.Item(8, 31).Value = 129369
$ws.Cells.Item(8, 32).Value = 0.22
$ws.Cells.Item(8, 33).Value = 761
$ws.Cells.Item(8, 34).Value = 2.69
$ws.Cells.Item(8, 35).Value = 44.98

# Row 9
$ws.Cells.Item(9, 4).Value = 204864
$ws.Cells.Item(9, 5).Value = 7433
$ws.Cells.Item(9, 7).Value = 4816
$ws.Cells.Item(9, 8).Value = 3529
$ws.Cells.Item(9, 9).Value = 3480
$ws.Cells.Item(9, 11).Value = 334927
$ws.Cells.Item(9, 12).Value = 159267
$ws.Cells.Item(9, 13).Value = 175660
$ws.Cells.Item(9, 14).Value = 172325
$ws.Cells.Item(9, 16).Value = 6670
$ws.Cells.Item(9, 17).Value = 19388
$ws.Cells.Item(9, 18).Value = -11907
$ws.Cells.Item(9, 19).Value = -5045
$ws.Cells.Item(9, 20).Value = 11618
$ws.Cells.Item(9, 21).Value = 7746
$ws.Cells.Item(9, 23).Value = 3.63
$ws.Cells.Item(9, 24).Value = 1.72
$ws.Cells.Item(9, 25).Value = 2.03
$ws.Cells.Item(9, 26).Value = 1.05
$ws.Cells.Item(9, 27).Value = 90.67
$ws.Cells.Item(9, 29).Value = 2608
$ws.Cells.Item(9, 30).Value = 10.85
$ws.Cells.Item(9, 31).Value = 131274
$ws.Cells.Item(9, 32).Value = 0.22
$ws.Cells.Item(9, 33).Value = 764
$ws.Cells.Item(9, 34).Value = 2.7
$ws.Cells.Item(9, 35).Value = 29.29
